$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1, G1) - copy formatting (bold, centered, bordered) from the
# existing header style used by E1, then set their text.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "2022_매출액"
$ws.Range("G1").Value = "2023_매출액"

# New data columns: 2022_매출액 (F) and 2023_매출액 (G) per industry row.
$data = @(
    @(2, 147, 215),
    @(3, 0, 0),
    @(4, 33328, 34799),
    @(5, 0, 0),
    @(6, 3708, 1468),
    @(7, 27960, 29088),
    @(8, 147478, 146413),
    @(9, 8843, 7205),
    @(10, 9875, 10871),
    @(11, 5640, 5381),
    @(12, 349, 177),
    @(13, 5394, 3905),
    @(14, 8611, 8752),
    @(15, 9713, 11168),
    @(16, 1776, 1966),
    @(17, 11505, 12586),
    @(18, 2408, 2856),
    @(19, 4511, 4939),
    @(20, 230012, 284892)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
}
